$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a NUMBER into a cell that is number-formatted as Text ("@").
# A plain ".Value = <number>" assignment on such a cell is stored as a text
# string (matches real Excel's text-cell coercion), but the target file
# keeps these as genuine numeric cells. Work around it by temporarily
# borrowing a "General" formatted cell's format (same border/fill/font),
# writing the number, then restoring the original Text format - all via
# Copy/PasteSpecial(Formats) so no new style entries are minted.
# ---------------------------------------------------------------------------
function Set-NumberInTextCell {
    param($ws, $row, $col, $value, $generalDonorAddr, $textDonorAddr)

    $ws.Range($generalDonorAddr).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Cells.Item($row, $col).Value = $value
    $ws.Range($textDonorAddr).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Update "Nb nouveaux cas positifs" (column C) for rows 636-643 ---
# (column B "Cumul cas positifs" is a running-total formula and recalculates
# automatically once C is updated, cascading through the rest of the table)
$ws.Cells.Item(636, 3).Value = 301
$ws.Cells.Item(637, 3).Value = 246
$ws.Cells.Item(638, 3).Value = 268
$ws.Cells.Item(639, 3).Value = 270
$ws.Cells.Item(640, 3).Value = 266
$ws.Cells.Item(641, 3).Value = 153
$ws.Cells.Item(642, 3).Value = 148
$ws.Cells.Item(643, 3).Value = 419

# --- Fill in the previously-empty daily rows 644-646 ---
# D650/C650 are untouched "General" formatted cells (same border as column
# D/C) used as temporary format donors; L650/M650 are untouched cells that
# still carry the original Text ("@") format for columns L/M.

# Row 644 : 2021-11-30
$ws.Cells.Item(644, 3).Value  = 315   # C - Nb nouveaux cas positifs
$ws.Cells.Item(644, 5).Value  = 10    # E - Nb nouvelles admissions a l'hopital
$ws.Cells.Item(644, 6).Value  = 4     # F - Patients COVID-19 intubes
$ws.Cells.Item(644, 7).Value  = 60    # G - Patients COVID-19 hospitalises hors SI
Set-NumberInTextCell $ws 644 12 0 "D650" "L650"   # L - Nb nouveaux deces a l'hopital
Set-NumberInTextCell $ws 644 13 0 "C650" "M650"   # M - Nb nouveaux deces extra-hospitaliers

# Row 645 : 2021-12-01
$ws.Cells.Item(645, 3).Value  = 239
$ws.Cells.Item(645, 5).Value  = 8
$ws.Cells.Item(645, 6).Value  = 4
$ws.Cells.Item(645, 7).Value  = 59
Set-NumberInTextCell $ws 645 12 2 "D650" "L650"
Set-NumberInTextCell $ws 645 13 0 "C650" "M650"

# Row 646 : 2021-12-02
$ws.Cells.Item(646, 3).Value  = 9
$ws.Cells.Item(646, 5).Value  = 8
$ws.Cells.Item(646, 6).Value  = 6
$ws.Cells.Item(646, 7).Value  = 60
Set-NumberInTextCell $ws 646 12 0 "D650" "L650"
Set-NumberInTextCell $ws 646 13 0 "C650" "M650"
